# Update metadata for the "municipio-nombre" column (column D).
# It used to be flagged as a measure (iaest-measure:municipio-nombre / medida / xsd:int);
# it is now recurated as a dimension, matching the pattern already used by the
# "provincia-nombre" column (E): sdmx-dimension:refArea / dim / URI-<entity>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
